$wb = $excel.ActiveWorkbook

# Update the root driver value on "updates_from_sub" sheet (B2),
# which feeds 'Calculations Rough Input'!D3 via VLOOKUP and cascades
# through the rest of the workbook's formulas.
$updates = $wb.Worksheets.Item("updates_from_sub")
$updates.Range("B2").Value = 10.72229629160968

# Update probe mass values (column F) on the "Probes" sheet from 210 to 212
$probes = $wb.Worksheets.Item("Probes")
for ($r = 2; $r -le 12; $r++) {
    $probes.Cells.Item($r, 6).Value = 212
}

# Make "Probes" the active sheet/tab (tabSelected moves there) and set its selection
$probes.Activate()
$probes.Range("H12").Select()

$wb.Save()
